# Apply crypto price/volume updates per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.315.09"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "2.589.77"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("D6").Value = "'143.86"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "2.599.63"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").Value = "'6.69"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("E11").Value = "  +3.44%  "
$ws.Range("D12").Value = "'0.158"
$ws.Range("E12").Value = "  +10.90%  "
$ws.Range("E13").Value = "  +2.59%  "
$ws.Range("D14").Value = "3.048.28"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "59.305.96"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").Value = "'22.57"
$ws.Range("E16").Value = "  +7.82%  "
$ws.Range("E17").Value = "  +3.78%  "
$ws.Range("D18").Value = "2.598.07"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").Value = "'336.65"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "'10.27"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("D22").Value = "'6.22"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'64.38"
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("D25").Value = "'0.457"
$ws.Range("E25").Value = "  +6.66%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("D29").Value = "0.0₃0783"
$ws.Range("E29").Value = "  +3.21%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("D33").Value = "'159.10"
$ws.Range("E33").Value = "  +3.22%  "
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("D37").Value = "'0.882"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").Value = "'0.877"
$ws.Range("E38").Value = "  -2.16%  "
$ws.Range("D39").Value = "'37.16"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("D41").Value = "'295.49"
$ws.Range("E41").Value = "  +4.09%  "
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  +2.80%  "
$ws.Range("D45").Value = "'0.594"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'19.36"
$ws.Range("E46").Value = "  +2.52%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0539"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").Value = "'125.82"
$ws.Range("E49").Value = "  +7.18%  "
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("D51").Value = "1.956.50"
$ws.Range("E51").Value = "  +0.42%  "
